# This script reproduces the "fix led locations and add to source code" edit:
# For ~107 data rows it fills in LED-location columns D (Latitude), E (Longitude),
# F (openLr Code) with placeholder 0s and a new computed column G (Free Flow Speed
# delta), and (for rows that had no data at all yet) marks column H ("Duplicate")
# as "Y". It also renumbers two LED Numbers (A92/A93) and updates the active
# selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1 - stamp the correct cell styles onto every target D/E/F cell by
# copy/pasting formats from template cells that already carry the right style
# (D92 -> style used by D/E i.e. Segoe UI font; A5 -> style used by F i.e.
# centered default font). This reuses the existing style indices instead of
# minting new ones.
# ---------------------------------------------------------------------------
$ws.Range("D92").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("E39").PasteSpecial(-4122)
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("E42").PasteSpecial(-4122)
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("E44").PasteSpecial(-4122)
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("E51").PasteSpecial(-4122)
$ws.Range("D53").PasteSpecial(-4122)
$ws.Range("E53").PasteSpecial(-4122)
$ws.Range("D54").PasteSpecial(-4122)
$ws.Range("E54").PasteSpecial(-4122)
$ws.Range("D55").PasteSpecial(-4122)
$ws.Range("E55").PasteSpecial(-4122)
$ws.Range("D56").PasteSpecial(-4122)
$ws.Range("E56").PasteSpecial(-4122)
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("E57").PasteSpecial(-4122)
$ws.Range("D58").PasteSpecial(-4122)
$ws.Range("E58").PasteSpecial(-4122)
$ws.Range("D59").PasteSpecial(-4122)
$ws.Range("E59").PasteSpecial(-4122)
$ws.Range("D61").PasteSpecial(-4122)
$ws.Range("E61").PasteSpecial(-4122)
$ws.Range("D62").PasteSpecial(-4122)
$ws.Range("E62").PasteSpecial(-4122)
$ws.Range("D63").PasteSpecial(-4122)
$ws.Range("E63").PasteSpecial(-4122)
$ws.Range("D64").PasteSpecial(-4122)
$ws.Range("E64").PasteSpecial(-4122)
$ws.Range("D65").PasteSpecial(-4122)
$ws.Range("E65").PasteSpecial(-4122)
$ws.Range("D66").PasteSpecial(-4122)
$ws.Range("E66").PasteSpecial(-4122)
$ws.Range("D70").PasteSpecial(-4122)
$ws.Range("E70").PasteSpecial(-4122)
$ws.Range("D76").PasteSpecial(-4122)
$ws.Range("E76").PasteSpecial(-4122)
$ws.Range("D94").PasteSpecial(-4122)
$ws.Range("E94").PasteSpecial(-4122)
$ws.Range("D99").PasteSpecial(-4122)
$ws.Range("E99").PasteSpecial(-4122)
$ws.Range("D101").PasteSpecial(-4122)
$ws.Range("E101").PasteSpecial(-4122)
$ws.Range("D110").PasteSpecial(-4122)
$ws.Range("E110").PasteSpecial(-4122)
$ws.Range("D125").PasteSpecial(-4122)
$ws.Range("E125").PasteSpecial(-4122)
$ws.Range("D127").PasteSpecial(-4122)
$ws.Range("E127").PasteSpecial(-4122)
$ws.Range("D132").PasteSpecial(-4122)
$ws.Range("E132").PasteSpecial(-4122)
$ws.Range("D133").PasteSpecial(-4122)
$ws.Range("E133").PasteSpecial(-4122)
$ws.Range("D134").PasteSpecial(-4122)
$ws.Range("E134").PasteSpecial(-4122)
$ws.Range("D160").PasteSpecial(-4122)
$ws.Range("E160").PasteSpecial(-4122)
$ws.Range("D161").PasteSpecial(-4122)
$ws.Range("E161").PasteSpecial(-4122)
$ws.Range("D225").PasteSpecial(-4122)
$ws.Range("E225").PasteSpecial(-4122)
$ws.Range("D226").PasteSpecial(-4122)
$ws.Range("E226").PasteSpecial(-4122)
$ws.Range("D227").PasteSpecial(-4122)
$ws.Range("E227").PasteSpecial(-4122)
$ws.Range("D228").PasteSpecial(-4122)
$ws.Range("E228").PasteSpecial(-4122)
$ws.Range("D229").PasteSpecial(-4122)
$ws.Range("E229").PasteSpecial(-4122)
$ws.Range("D231").PasteSpecial(-4122)
$ws.Range("E231").PasteSpecial(-4122)
$ws.Range("D232").PasteSpecial(-4122)
$ws.Range("E232").PasteSpecial(-4122)
$ws.Range("D233").PasteSpecial(-4122)
$ws.Range("E233").PasteSpecial(-4122)
$ws.Range("D234").PasteSpecial(-4122)
$ws.Range("E234").PasteSpecial(-4122)
$ws.Range("D235").PasteSpecial(-4122)
$ws.Range("E235").PasteSpecial(-4122)
$ws.Range("D237").PasteSpecial(-4122)
$ws.Range("E237").PasteSpecial(-4122)
$ws.Range("D239").PasteSpecial(-4122)
$ws.Range("E239").PasteSpecial(-4122)
$ws.Range("D241").PasteSpecial(-4122)
$ws.Range("E241").PasteSpecial(-4122)
$ws.Range("D243").PasteSpecial(-4122)
$ws.Range("E243").PasteSpecial(-4122)
$ws.Range("D245").PasteSpecial(-4122)
$ws.Range("E245").PasteSpecial(-4122)
$ws.Range("D246").PasteSpecial(-4122)
$ws.Range("E246").PasteSpecial(-4122)
$ws.Range("D247").PasteSpecial(-4122)
$ws.Range("E247").PasteSpecial(-4122)
$ws.Range("D249").PasteSpecial(-4122)
$ws.Range("E249").PasteSpecial(-4122)
$ws.Range("D250").PasteSpecial(-4122)
$ws.Range("E250").PasteSpecial(-4122)
$ws.Range("D251").PasteSpecial(-4122)
$ws.Range("E251").PasteSpecial(-4122)
$ws.Range("D255").PasteSpecial(-4122)
$ws.Range("E255").PasteSpecial(-4122)
$ws.Range("D257").PasteSpecial(-4122)
$ws.Range("E257").PasteSpecial(-4122)
$ws.Range("D259").PasteSpecial(-4122)
$ws.Range("E259").PasteSpecial(-4122)
$ws.Range("D261").PasteSpecial(-4122)
$ws.Range("E261").PasteSpecial(-4122)
$ws.Range("D263").PasteSpecial(-4122)
$ws.Range("E263").PasteSpecial(-4122)
$ws.Range("D264").PasteSpecial(-4122)
$ws.Range("E264").PasteSpecial(-4122)
$ws.Range("D266").PasteSpecial(-4122)
$ws.Range("E266").PasteSpecial(-4122)
$ws.Range("D267").PasteSpecial(-4122)
$ws.Range("E267").PasteSpecial(-4122)
$ws.Range("D269").PasteSpecial(-4122)
$ws.Range("E269").PasteSpecial(-4122)
$ws.Range("D270").PasteSpecial(-4122)
$ws.Range("E270").PasteSpecial(-4122)
$ws.Range("D271").PasteSpecial(-4122)
$ws.Range("E271").PasteSpecial(-4122)
$ws.Range("D273").PasteSpecial(-4122)
$ws.Range("E273").PasteSpecial(-4122)
$ws.Range("D274").PasteSpecial(-4122)
$ws.Range("E274").PasteSpecial(-4122)
$ws.Range("D289").PasteSpecial(-4122)
$ws.Range("E289").PasteSpecial(-4122)
$ws.Range("D319").PasteSpecial(-4122)
$ws.Range("E319").PasteSpecial(-4122)
$ws.Range("D321").PasteSpecial(-4122)
$ws.Range("E321").PasteSpecial(-4122)
$ws.Range("D335").PasteSpecial(-4122)
$ws.Range("E335").PasteSpecial(-4122)
$ws.Range("D337").PasteSpecial(-4122)
$ws.Range("E337").PasteSpecial(-4122)
$ws.Range("D339").PasteSpecial(-4122)
$ws.Range("E339").PasteSpecial(-4122)
$ws.Range("D341").PasteSpecial(-4122)
$ws.Range("E341").PasteSpecial(-4122)
$ws.Range("D359").PasteSpecial(-4122)
$ws.Range("E359").PasteSpecial(-4122)
$ws.Range("D373").PasteSpecial(-4122)
$ws.Range("E373").PasteSpecial(-4122)
$ws.Range("D374").PasteSpecial(-4122)
$ws.Range("E374").PasteSpecial(-4122)
$ws.Range("D389").PasteSpecial(-4122)
$ws.Range("E389").PasteSpecial(-4122)
$ws.Range("D391").PasteSpecial(-4122)
$ws.Range("E391").PasteSpecial(-4122)
$ws.Range("D401").PasteSpecial(-4122)
$ws.Range("E401").PasteSpecial(-4122)
$ws.Range("D419").PasteSpecial(-4122)
$ws.Range("E419").PasteSpecial(-4122)
$ws.Range("D425").PasteSpecial(-4122)
$ws.Range("E425").PasteSpecial(-4122)
$ws.Range("D428").PasteSpecial(-4122)
$ws.Range("E428").PasteSpecial(-4122)
$ws.Range("D429").PasteSpecial(-4122)
$ws.Range("E429").PasteSpecial(-4122)
$ws.Range("D430").PasteSpecial(-4122)
$ws.Range("E430").PasteSpecial(-4122)
$ws.Range("D432").PasteSpecial(-4122)
$ws.Range("E432").PasteSpecial(-4122)
$ws.Range("D440").PasteSpecial(-4122)
$ws.Range("E440").PasteSpecial(-4122)
$ws.Range("D473").PasteSpecial(-4122)
$ws.Range("E473").PasteSpecial(-4122)
$ws.Range("D474").PasteSpecial(-4122)
$ws.Range("E474").PasteSpecial(-4122)
$ws.Range("D475").PasteSpecial(-4122)
$ws.Range("E475").PasteSpecial(-4122)
$ws.Range("D477").PasteSpecial(-4122)
$ws.Range("E477").PasteSpecial(-4122)
$ws.Range("D478").PasteSpecial(-4122)
$ws.Range("E478").PasteSpecial(-4122)
$ws.Range("D479").PasteSpecial(-4122)
$ws.Range("E479").PasteSpecial(-4122)
$ws.Range("D482").PasteSpecial(-4122)
$ws.Range("E482").PasteSpecial(-4122)
$ws.Range("D483").PasteSpecial(-4122)
$ws.Range("E483").PasteSpecial(-4122)
$ws.Range("D497").PasteSpecial(-4122)
$ws.Range("E497").PasteSpecial(-4122)
$ws.Range("D499").PasteSpecial(-4122)
$ws.Range("E499").PasteSpecial(-4122)
$ws.Range("D512").PasteSpecial(-4122)
$ws.Range("E512").PasteSpecial(-4122)
$ws.Range("D518").PasteSpecial(-4122)
$ws.Range("E518").PasteSpecial(-4122)
$ws.Range("D551").PasteSpecial(-4122)
$ws.Range("E551").PasteSpecial(-4122)
$ws.Range("D596").PasteSpecial(-4122)
$ws.Range("E596").PasteSpecial(-4122)
$ws.Range("D601").PasteSpecial(-4122)
$ws.Range("E601").PasteSpecial(-4122)
$ws.Range("D628").PasteSpecial(-4122)
$ws.Range("E628").PasteSpecial(-4122)
$ws.Range("D641").PasteSpecial(-4122)
$ws.Range("E641").PasteSpecial(-4122)
$ws.Range("D645").PasteSpecial(-4122)
$ws.Range("E645").PasteSpecial(-4122)
$ws.Range("D647").PasteSpecial(-4122)
$ws.Range("E647").PasteSpecial(-4122)
$ws.Range("D648").PasteSpecial(-4122)
$ws.Range("E648").PasteSpecial(-4122)

$ws.Range("A5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("F42").PasteSpecial(-4122)
$ws.Range("F44").PasteSpecial(-4122)
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("F61").PasteSpecial(-4122)
$ws.Range("F62").PasteSpecial(-4122)
$ws.Range("F63").PasteSpecial(-4122)
$ws.Range("F64").PasteSpecial(-4122)
$ws.Range("F65").PasteSpecial(-4122)
$ws.Range("F66").PasteSpecial(-4122)
$ws.Range("F70").PasteSpecial(-4122)
$ws.Range("F76").PasteSpecial(-4122)
$ws.Range("F94").PasteSpecial(-4122)
$ws.Range("F99").PasteSpecial(-4122)
$ws.Range("F101").PasteSpecial(-4122)
$ws.Range("F110").PasteSpecial(-4122)
$ws.Range("F125").PasteSpecial(-4122)
$ws.Range("F127").PasteSpecial(-4122)
$ws.Range("F132").PasteSpecial(-4122)
$ws.Range("F133").PasteSpecial(-4122)
$ws.Range("F134").PasteSpecial(-4122)
$ws.Range("F160").PasteSpecial(-4122)
$ws.Range("F161").PasteSpecial(-4122)
$ws.Range("F225").PasteSpecial(-4122)
$ws.Range("F226").PasteSpecial(-4122)
$ws.Range("F227").PasteSpecial(-4122)
$ws.Range("F228").PasteSpecial(-4122)
$ws.Range("F229").PasteSpecial(-4122)
$ws.Range("F231").PasteSpecial(-4122)
$ws.Range("F232").PasteSpecial(-4122)
$ws.Range("F233").PasteSpecial(-4122)
$ws.Range("F234").PasteSpecial(-4122)
$ws.Range("F235").PasteSpecial(-4122)
$ws.Range("F237").PasteSpecial(-4122)
$ws.Range("F239").PasteSpecial(-4122)
$ws.Range("F241").PasteSpecial(-4122)
$ws.Range("F243").PasteSpecial(-4122)
$ws.Range("F245").PasteSpecial(-4122)
$ws.Range("F246").PasteSpecial(-4122)
$ws.Range("F247").PasteSpecial(-4122)
$ws.Range("F249").PasteSpecial(-4122)
$ws.Range("F250").PasteSpecial(-4122)
$ws.Range("F251").PasteSpecial(-4122)
$ws.Range("F255").PasteSpecial(-4122)
$ws.Range("F257").PasteSpecial(-4122)
$ws.Range("F259").PasteSpecial(-4122)
$ws.Range("F261").PasteSpecial(-4122)
$ws.Range("F263").PasteSpecial(-4122)
$ws.Range("F264").PasteSpecial(-4122)
$ws.Range("F266").PasteSpecial(-4122)
$ws.Range("F267").PasteSpecial(-4122)
$ws.Range("F269").PasteSpecial(-4122)
$ws.Range("F270").PasteSpecial(-4122)
$ws.Range("F271").PasteSpecial(-4122)
$ws.Range("F273").PasteSpecial(-4122)
$ws.Range("F274").PasteSpecial(-4122)
$ws.Range("F289").PasteSpecial(-4122)
$ws.Range("F319").PasteSpecial(-4122)
$ws.Range("F321").PasteSpecial(-4122)
$ws.Range("F335").PasteSpecial(-4122)
$ws.Range("F337").PasteSpecial(-4122)
$ws.Range("F339").PasteSpecial(-4122)
$ws.Range("F341").PasteSpecial(-4122)
$ws.Range("F359").PasteSpecial(-4122)
$ws.Range("F373").PasteSpecial(-4122)
$ws.Range("F374").PasteSpecial(-4122)
$ws.Range("F389").PasteSpecial(-4122)
$ws.Range("F391").PasteSpecial(-4122)
$ws.Range("F401").PasteSpecial(-4122)
$ws.Range("F419").PasteSpecial(-4122)
$ws.Range("F425").PasteSpecial(-4122)
$ws.Range("F428").PasteSpecial(-4122)
$ws.Range("F429").PasteSpecial(-4122)
$ws.Range("F430").PasteSpecial(-4122)
$ws.Range("F432").PasteSpecial(-4122)
$ws.Range("F440").PasteSpecial(-4122)
$ws.Range("F473").PasteSpecial(-4122)
$ws.Range("F474").PasteSpecial(-4122)
$ws.Range("F475").PasteSpecial(-4122)
$ws.Range("F477").PasteSpecial(-4122)
$ws.Range("F478").PasteSpecial(-4122)
$ws.Range("F479").PasteSpecial(-4122)
$ws.Range("F482").PasteSpecial(-4122)
$ws.Range("F483").PasteSpecial(-4122)
$ws.Range("F497").PasteSpecial(-4122)
$ws.Range("F499").PasteSpecial(-4122)
$ws.Range("F512").PasteSpecial(-4122)
$ws.Range("F518").PasteSpecial(-4122)
$ws.Range("F551").PasteSpecial(-4122)
$ws.Range("F596").PasteSpecial(-4122)
$ws.Range("F601").PasteSpecial(-4122)
$ws.Range("F628").PasteSpecial(-4122)
$ws.Range("F641").PasteSpecial(-4122)
$ws.Range("F645").PasteSpecial(-4122)
$ws.Range("F647").PasteSpecial(-4122)
$ws.Range("F648").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2 - write the actual values: D/E/F become 0, G becomes the row-specific
# delta that was computed for that LED.
# ---------------------------------------------------------------------------
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = -2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = -2
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = -4
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = -4
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = -19
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = -20
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = -21
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = -26
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = -25
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = -25
$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = -25
$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = -25
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = -25
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = -25
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = -25
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = -27
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = -27
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = -27
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = -27
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = -27
$ws.Range("D66").Value = 0
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = -27
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = -34
$ws.Range("D76").Value = 0
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = -37
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 0
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = -73
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = -76
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = -76
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 0
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = -81
$ws.Range("D125").Value = 0
$ws.Range("E125").Value = 0
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = -107
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 0
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = -107
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 0
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = -146
$ws.Range("D133").Value = 0
$ws.Range("E133").Value = 0
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = -146
$ws.Range("D134").Value = 0
$ws.Range("E134").Value = 0
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = -146
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 0
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = -250
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 0
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = -250
$ws.Range("D225").Value = 0
$ws.Range("E225").Value = 0
$ws.Range("F225").Value = 0
$ws.Range("G225").Value = -61
$ws.Range("D226").Value = 0
$ws.Range("E226").Value = 0
$ws.Range("F226").Value = 0
$ws.Range("G226").Value = -61
$ws.Range("D227").Value = 0
$ws.Range("E227").Value = 0
$ws.Range("F227").Value = 0
$ws.Range("G227").Value = -61
$ws.Range("D228").Value = 0
$ws.Range("E228").Value = 0
$ws.Range("F228").Value = 0
$ws.Range("G228").Value = -61
$ws.Range("D229").Value = 0
$ws.Range("E229").Value = 0
$ws.Range("F229").Value = 0
$ws.Range("G229").Value = -61
$ws.Range("D231").Value = 0
$ws.Range("E231").Value = 0
$ws.Range("F231").Value = 0
$ws.Range("G231").Value = -61
$ws.Range("D232").Value = 0
$ws.Range("E232").Value = 0
$ws.Range("F232").Value = 0
$ws.Range("G232").Value = -61
$ws.Range("D233").Value = 0
$ws.Range("E233").Value = 0
$ws.Range("F233").Value = 0
$ws.Range("G233").Value = -61
$ws.Range("D234").Value = 0
$ws.Range("E234").Value = 0
$ws.Range("F234").Value = 0
$ws.Range("G234").Value = -61
$ws.Range("D235").Value = 0
$ws.Range("E235").Value = 0
$ws.Range("F235").Value = 0
$ws.Range("G235").Value = -61
$ws.Range("D237").Value = 0
$ws.Range("E237").Value = 0
$ws.Range("F237").Value = 0
$ws.Range("G237").Value = -85
$ws.Range("D239").Value = 0
$ws.Range("E239").Value = 0
$ws.Range("F239").Value = 0
$ws.Range("G239").Value = -85
$ws.Range("D241").Value = 0
$ws.Range("E241").Value = 0
$ws.Range("F241").Value = 0
$ws.Range("G241").Value = -87
$ws.Range("D243").Value = 0
$ws.Range("E243").Value = 0
$ws.Range("F243").Value = 0
$ws.Range("G243").Value = -87
$ws.Range("D245").Value = 0
$ws.Range("E245").Value = 0
$ws.Range("F245").Value = 0
$ws.Range("G245").Value = -89
$ws.Range("D246").Value = 0
$ws.Range("E246").Value = 0
$ws.Range("F246").Value = 0
$ws.Range("G246").Value = -89
$ws.Range("D247").Value = 0
$ws.Range("E247").Value = 0
$ws.Range("F247").Value = 0
$ws.Range("G247").Value = -89
$ws.Range("D249").Value = 0
$ws.Range("E249").Value = 0
$ws.Range("F249").Value = 0
$ws.Range("G249").Value = -88
$ws.Range("D250").Value = 0
$ws.Range("E250").Value = 0
$ws.Range("F250").Value = 0
$ws.Range("G250").Value = -88
$ws.Range("D251").Value = 0
$ws.Range("E251").Value = 0
$ws.Range("F251").Value = 0
$ws.Range("G251").Value = -88
$ws.Range("D255").Value = 0
$ws.Range("E255").Value = 0
$ws.Range("F255").Value = 0
$ws.Range("G255").Value = -112
$ws.Range("D257").Value = 0
$ws.Range("E257").Value = 0
$ws.Range("F257").Value = 0
$ws.Range("G257").Value = -112
$ws.Range("D259").Value = 0
$ws.Range("E259").Value = 0
$ws.Range("F259").Value = 0
$ws.Range("G259").Value = -114
$ws.Range("D261").Value = 0
$ws.Range("E261").Value = 0
$ws.Range("F261").Value = 0
$ws.Range("G261").Value = -114
$ws.Range("D263").Value = 0
$ws.Range("E263").Value = 0
$ws.Range("F263").Value = 0
$ws.Range("G263").Value = -116
$ws.Range("D264").Value = 0
$ws.Range("E264").Value = 0
$ws.Range("F264").Value = 0
$ws.Range("G264").Value = -116
$ws.Range("D266").Value = 0
$ws.Range("E266").Value = 0
$ws.Range("F266").Value = 0
$ws.Range("G266").Value = -116
$ws.Range("D267").Value = 0
$ws.Range("E267").Value = 0
$ws.Range("F267").Value = 0
$ws.Range("G267").Value = -116
$ws.Range("D269").Value = 0
$ws.Range("E269").Value = 0
$ws.Range("F269").Value = 0
$ws.Range("G269").Value = -155
$ws.Range("D270").Value = 0
$ws.Range("E270").Value = 0
$ws.Range("F270").Value = 0
$ws.Range("G270").Value = -155
$ws.Range("D271").Value = 0
$ws.Range("E271").Value = 0
$ws.Range("F271").Value = 0
$ws.Range("G271").Value = -155
$ws.Range("D273").Value = 0
$ws.Range("E273").Value = 0
$ws.Range("F273").Value = 0
$ws.Range("G273").Value = -155
$ws.Range("D274").Value = 0
$ws.Range("E274").Value = 0
$ws.Range("F274").Value = 0
$ws.Range("G274").Value = -155
$ws.Range("D289").Value = 0
$ws.Range("E289").Value = 0
$ws.Range("F289").Value = 0
$ws.Range("G289").Value = -192
$ws.Range("D319").Value = 0
$ws.Range("E319").Value = 0
$ws.Range("F319").Value = 0
$ws.Range("G319").Value = -207
$ws.Range("D321").Value = 0
$ws.Range("E321").Value = 0
$ws.Range("F321").Value = 0
$ws.Range("G321").Value = -207
$ws.Range("D335").Value = 0
$ws.Range("E335").Value = 0
$ws.Range("F335").Value = 0
$ws.Range("G335").Value = -233
$ws.Range("D337").Value = 0
$ws.Range("E337").Value = 0
$ws.Range("F337").Value = 0
$ws.Range("G337").Value = -233
$ws.Range("D339").Value = 0
$ws.Range("E339").Value = 0
$ws.Range("F339").Value = 0
$ws.Range("G339").Value = -329
$ws.Range("D341").Value = 0
$ws.Range("E341").Value = 0
$ws.Range("F341").Value = 0
$ws.Range("G341").Value = -329
$ws.Range("D359").Value = 0
$ws.Range("E359").Value = 0
$ws.Range("F359").Value = 0
$ws.Range("G359").Value = -54
$ws.Range("D373").Value = 0
$ws.Range("E373").Value = 0
$ws.Range("F373").Value = 0
$ws.Range("G373").Value = -70
$ws.Range("D374").Value = 0
$ws.Range("E374").Value = 0
$ws.Range("F374").Value = 0
$ws.Range("G374").Value = -70
$ws.Range("D389").Value = 0
$ws.Range("E389").Value = 0
$ws.Range("F389").Value = 0
$ws.Range("G389").Value = -96
$ws.Range("D391").Value = 0
$ws.Range("E391").Value = 0
$ws.Range("F391").Value = 0
$ws.Range("G391").Value = -96
$ws.Range("D401").Value = 0
$ws.Range("E401").Value = 0
$ws.Range("F401").Value = 0
$ws.Range("G401").Value = -138
$ws.Range("D419").Value = 0
$ws.Range("E419").Value = 0
$ws.Range("F419").Value = 0
$ws.Range("G419").Value = -166
$ws.Range("D425").Value = 0
$ws.Range("E425").Value = 0
$ws.Range("F425").Value = 0
$ws.Range("G425").Value = -169
$ws.Range("D428").Value = 0
$ws.Range("E428").Value = 0
$ws.Range("F428").Value = 0
$ws.Range("G428").Value = -170
$ws.Range("D429").Value = 0
$ws.Range("E429").Value = 0
$ws.Range("F429").Value = 0
$ws.Range("G429").Value = -170
$ws.Range("D430").Value = 0
$ws.Range("E430").Value = 0
$ws.Range("F430").Value = 0
$ws.Range("G430").Value = -170
$ws.Range("D432").Value = 0
$ws.Range("E432").Value = 0
$ws.Range("F432").Value = 0
$ws.Range("G432").Value = -171
$ws.Range("D440").Value = 0
$ws.Range("E440").Value = 0
$ws.Range("F440").Value = 0
$ws.Range("G440").Value = -257
$ws.Range("D473").Value = 0
$ws.Range("E473").Value = 0
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = -118
$ws.Range("D474").Value = 0
$ws.Range("E474").Value = 0
$ws.Range("F474").Value = 0
$ws.Range("G474").Value = -118
$ws.Range("D475").Value = 0
$ws.Range("E475").Value = 0
$ws.Range("F475").Value = 0
$ws.Range("G475").Value = -118
$ws.Range("D477").Value = 0
$ws.Range("E477").Value = 0
$ws.Range("F477").Value = 0
$ws.Range("G477").Value = -118
$ws.Range("D478").Value = 0
$ws.Range("E478").Value = 0
$ws.Range("F478").Value = 0
$ws.Range("G478").Value = -118
$ws.Range("D479").Value = 0
$ws.Range("E479").Value = 0
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = -118
$ws.Range("D482").Value = 0
$ws.Range("E482").Value = 0
$ws.Range("F482").Value = 0
$ws.Range("G482").Value = -122
$ws.Range("D483").Value = 0
$ws.Range("E483").Value = 0
$ws.Range("F483").Value = 0
$ws.Range("G483").Value = -122
$ws.Range("D497").Value = 0
$ws.Range("E497").Value = 0
$ws.Range("F497").Value = 0
$ws.Range("G497").Value = -131
$ws.Range("D499").Value = 0
$ws.Range("E499").Value = 0
$ws.Range("F499").Value = 0
$ws.Range("G499").Value = -131
$ws.Range("D512").Value = 0
$ws.Range("E512").Value = 0
$ws.Range("F512").Value = 0
$ws.Range("G512").Value = -235
$ws.Range("D518").Value = 0
$ws.Range("E518").Value = 0
$ws.Range("F518").Value = 0
$ws.Range("G518").Value = -238
$ws.Range("D551").Value = 0
$ws.Range("E551").Value = 0
$ws.Range("F551").Value = 0
$ws.Range("G551").Value = -291
$ws.Range("D596").Value = 0
$ws.Range("E596").Value = 0
$ws.Range("F596").Value = 0
$ws.Range("G596").Value = -185
$ws.Range("D601").Value = 0
$ws.Range("E601").Value = 0
$ws.Range("F601").Value = 0
$ws.Range("G601").Value = -188
$ws.Range("D628").Value = 0
$ws.Range("E628").Value = 0
$ws.Range("F628").Value = 0
$ws.Range("G628").Value = -219
$ws.Range("D641").Value = 0
$ws.Range("E641").Value = 0
$ws.Range("F641").Value = 0
$ws.Range("G641").Value = -320
$ws.Range("D645").Value = 0
$ws.Range("E645").Value = 0
$ws.Range("F645").Value = 0
$ws.Range("G645").Value = -322
$ws.Range("D647").Value = 0
$ws.Range("E647").Value = 0
$ws.Range("F647").Value = 0
$ws.Range("G647").Value = -322
$ws.Range("D648").Value = 0
$ws.Range("E648").Value = 0
$ws.Range("F648").Value = 0
$ws.Range("G648").Value = -322

# ---------------------------------------------------------------------------
# Step 3 - these rows had no data at all past column C before; mark their new
# "Duplicate" column H as "Y" too.
# ---------------------------------------------------------------------------
$ws.Range("H161").Value = "Y"
$ws.Range("H231").Value = "Y"
$ws.Range("H232").Value = "Y"
$ws.Range("H233").Value = "Y"
$ws.Range("H234").Value = "Y"
$ws.Range("H235").Value = "Y"
$ws.Range("H249").Value = "Y"
$ws.Range("H250").Value = "Y"
$ws.Range("H251").Value = "Y"

# ---------------------------------------------------------------------------
# Step 4 - renumber the LED Number for the pair of rows that moved from #46 to
# #73 in the source data.
# ---------------------------------------------------------------------------
$ws.Range("A92").Value = 73
$ws.Range("A93").Value = 73

# ---------------------------------------------------------------------------
# Step 5 - move the active selection (the workbook was left scrolled/selected
# at a different cell before saving).
# ---------------------------------------------------------------------------
$ws.Range("H18").Select()
